$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.764364242553711
$ws.Range("B1").Value = 2.531865358352661
$ws.Range("C1").Value = 1.903707504272461
$ws.Range("D1").Value = 1.798515558242798
$ws.Range("E1").Value = 1.896444082260132
